$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("挑战组")
$ws1.Range("H1").Value = "link"
$ws1.Range("H1").Font.Bold = $true
$ws1.Range("H1").HorizontalAlignment = -4108
$ws1.Range("H1").VerticalAlignment = -4160
$ws1.Range("H1").Borders.Item(7).LineStyle = 1
$ws1.Range("H1").Borders.Item(7).ColorIndex = -4105
$ws1.Range("H1").Borders.Item(10).LineStyle = 1
$ws1.Range("H1").Borders.Item(10).ColorIndex = -4105
